$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")
$lo = $ws.ListObjects.Item(1)

# --- Expand the table to make room for the new "Unit" column (54 cols x 4 rows) ---
$lo.Resize($ws.Range("A1:BB4"))

# --- Insert a blank column at AA, shifting "altitude"-and-later columns one to the right ---
$ws.Range("AA1:AA4").Insert(-4161)

# --- Re-write every header cell so the table column names re-sync with the shifted text ---
$ws.Range("A1").Value = "Input [Source Name]"
$ws.Range("B1").Value = "Characteristic [MIAPPE version]"
$ws.Range("C1").Value = "Term Source REF (MIAPPE:0007)"
$ws.Range("D1").Value = "Term Accession Number (MIAPPE:0007)"
$ws.Range("E1").Value = "Characteristic [Start date of study]"
$ws.Range("F1").Value = "Term Source REF (MIAPPE:0013)"
$ws.Range("G1").Value = "Term Accession Number (MIAPPE:0013)"
$ws.Range("H1").Value = "Characteristic [End date of study]"
$ws.Range("I1").Value = "Term Source REF (MIAPPE:0014)"
$ws.Range("J1").Value = "Term Accession Number (MIAPPE:0014)"
$ws.Range("K1").Value = "Characteristic [Contact institution]"
$ws.Range("L1").Value = "Term Source REF (MIAPPE:0015)"
$ws.Range("M1").Value = "Term Accession Number (MIAPPE:0015)"
$ws.Range("N1").Value = "Characteristic [Geographic location (country)]"
$ws.Range("O1").Value = "Term Source REF (MIAPPE:0016)"
$ws.Range("P1").Value = "Term Accession Number (MIAPPE:0016)"
$ws.Range("Q1").Value = "Characteristic [Experimental site name]"
$ws.Range("R1").Value = "Term Source REF (MIAPPE:0017)"
$ws.Range("S1").Value = "Term Accession Number (MIAPPE:0017)"
$ws.Range("T1").Value = "Characteristic [Geographic location (latitude)]"
$ws.Range("U1").Value = "Term Source REF (MIAPPE:0018)"
$ws.Range("V1").Value = "Term Accession Number (MIAPPE:0018)"
$ws.Range("W1").Value = "Characteristic [Geographic location (longitude)]"
$ws.Range("X1").Value = "Term Source REF (MIAPPE:0019)"
$ws.Range("Y1").Value = "Term Accession Number (MIAPPE:0019)"
$ws.Range("Z1").Value = "Characteristic [Geographic location (altitude)]"
$ws.Range("AA1").Value = "Unit"
$ws.Range("AB1").Value = "Term Source REF (MIAPPE:0020)"
$ws.Range("AC1").Value = "Term Accession Number (MIAPPE:0020)"
$ws.Range("AD1").Value = "Characteristic [Description of the experimental design]"
$ws.Range("AE1").Value = "Term Source REF (MIAPPE:0021)"
$ws.Range("AF1").Value = "Term Accession Number (MIAPPE:0021)"
$ws.Range("AG1").Value = "Characteristic [Type of experimental design]"
$ws.Range("AH1").Value = "Term Source REF (MIAPPE:0022)"
$ws.Range("AI1").Value = "Term Accession Number (MIAPPE:0022)"
$ws.Range("AJ1").Value = "Characteristic [Observation unit level hierarchy]"
$ws.Range("AK1").Value = "Term Source REF (MIAPPE:0023)"
$ws.Range("AL1").Value = "Term Accession Number (MIAPPE:0023)"
$ws.Range("AM1").Value = "Characteristic [Observation unit description]"
$ws.Range("AN1").Value = "Term Source REF (MIAPPE:0024)"
$ws.Range("AO1").Value = "Term Accession Number (MIAPPE:0024)"
$ws.Range("AP1").Value = "Characteristic [Description of growth facility]"
$ws.Range("AQ1").Value = "Term Source REF (MIAPPE:0025)"
$ws.Range("AR1").Value = "Term Accession Number (MIAPPE:0025)"
$ws.Range("AS1").Value = "Characteristic [Type of growth facility]"
$ws.Range("AT1").Value = "Term Source REF (MIAPPE:0026)"
$ws.Range("AU1").Value = "Term Accession Number (MIAPPE:0026)"
$ws.Range("AV1").Value = "Characteristic [Cultural practices]"
$ws.Range("AW1").Value = "Term Source REF (MIAPPE:0027)"
$ws.Range("AX1").Value = "Term Accession Number (MIAPPE:0027)"
$ws.Range("AY1").Value = "Characteristic [Map of experimental design]"
$ws.Range("AZ1").Value = "Term Source REF (MIAPPE:0028)"
$ws.Range("BA1").Value = "Term Accession Number (MIAPPE:0028)"
$ws.Range("BB1").Value = "Output [Sample Name]"

# --- Remove the now-superfluous blank example rows (3 and 4); one example row remains ---
$ws.Rows("4:4").Delete()
$ws.Rows("3:3").Delete()

# --- Force plain-text storage for values Excel would otherwise reinterpret as number/date ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("T2").NumberFormat = "@"
$ws.Range("W2").NumberFormat = "@"
$ws.Range("Z2").NumberFormat = "@"

# --- Populate row 2 with the MIAPPE example values ---
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "1.1"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "2002-04-04"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "2002-11-27"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "CNRS, Gif-sur-Yvette, France"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = "France"
$ws.Range("O2").Value = "NCIT"
$ws.Range("P2").Value = "http://purl.obolibrary.org/obo/NCIT_C16592"
$ws.Range("Q2").Value = " Domaine expérimental de Melgueil - 34130 Mauguio - France"
$ws.Range("R2").Value = ""
$ws.Range("S2").Value = ""
$ws.Range("T2").Value = "+43.619264"
$ws.Range("U2").Value = ""
$ws.Range("V2").Value = ""
$ws.Range("W2").Value = "+3.967454"
$ws.Range("X2").Value = ""
$ws.Range("Y2").Value = ""
$ws.Range("Z2").Value = "100"
$ws.Range("AA2").Value = "meter"
$ws.Range("AB2").Value = "UO"
$ws.Range("AC2").Value = "http://purl.obolibrary.org/obo/UO_0000008"
$ws.Range("AD2").Value = "Lines were repeated twice at each location using a complete block design. In order to limit competition effects, each block was organized into four sub-blocks corresponding to earliness groups based on a priori information. "
$ws.Range("AE2").Value = ""
$ws.Range("AF2").Value = ""
$ws.Range("AG2").Value = "randomized complete block design"
$ws.Range("AH2").Value = "OBI"
$ws.Range("AI2").Value = "http://purl.obolibrary.org/obo/OBI_0500007"
$ws.Range("AJ2").Value = "block>rep>plot"
$ws.Range("AK2").Value = ""
$ws.Range("AL2").Value = ""
$ws.Range("AM2").Value = "Observation units consisted in individual plots themselves consisting of a row of 15 plants at a density of approximately six plants per square meter."
$ws.Range("AN2").Value = ""
$ws.Range("AO2").Value = ""
$ws.Range("AP2").Value = ""
$ws.Range("AQ2").Value = ""
$ws.Range("AR2").Value = ""
$ws.Range("AS2").Value = "field environment condition"
$ws.Range("AT2").Value = ""
$ws.Range("AU2").Value = ""
$ws.Range("AV2").Value = "Irrigation was applied according needs during summer to prevent water stress."
$ws.Range("AW2").Value = ""
$ws.Range("AX2").Value = ""
$ws.Range("AY2").Value = "https://urgi.versailles.inra.fr/files/ephesis/181000503/181000503_plan.xls"
$ws.Range("AZ2").Value = ""
$ws.Range("BA2").Value = ""
$ws.Range("BB2").Value = ""

Write-Output "done"
